$wb = $excel.ActiveWorkbook

# Remember how many sheets existed originally so we can remove them once the
# replacement sheets are in place (a workbook can never have 0 sheets).
$origCount = $wb.Worksheets.Count
$lastOrig = $wb.Worksheets.Item($origCount)

function New-TaskOrderSheet($sheetName, $values, $after) {
    $ws = $wb.Worksheets.Add($null, $after)
    $ws.Name = $sheetName

    $ws.Range("B1").Value = "task_order"

    $r = 2
    for ($i = 0; $i -lt $values.Count; $i++) {
        $ws.Cells.Item($r, 1).Value = $i
        $ws.Cells.Item($r, 2).Value = $values[$i]
        $r++
    }

    # Format the header + index column together as a single union range so the
    # engine reuses one shared (bold / thin-border / center-top) style, matching
    # the original sheets instead of minting a new font/xf per cell touched.
    $lastRow = $values.Count + 1
    $styled = $ws.Range("B1,A2:A$lastRow")
    $styled.Font.Bold = $true
    $styled.HorizontalAlignment = -4108
    $styled.VerticalAlignment = -4160
    $styled.Borders.LineStyle = 1

    return $ws
}

$s1 = New-TaskOrderSheet "NB_TO-1651589019224096" @("TB-16515890188581567.csv", "OB-1651589017942814.csv", "ZB-match_9-16515890174383495.csv", "TB-16515890191928487.csv", "OB-16515890177531264.csv", "ZB-match_3-16515890174539745.csv", "ZB-match_3-16515890172065299.csv", "OB-16515890182106678.csv", "TB-16515890187622979.csv") $lastOrig
$s2 = New-TaskOrderSheet "TOL_TO-16515890192709715" @("MM_stims-16515890192397597.csv", "ZM_stims-1651589019224096.csv", "MM_stims-1651589019255347.csv", "ZM_stims-16515890192397597.csv", "MM_stims-16515890192709715.csv", "ZM_stims-1651589019255347.csv") $s1
$s3 = New-TaskOrderSheet "vSAT_TO-16515890193334723" @("SAT_stims-16515890192709715.csv", "SAT_stims-16515890192865949.csv", "vSAT_stims-16515890193178465.csv", "vSAT_stims-1651589019302222.csv") $s2
$s4 = New-TaskOrderSheet "RS_TO-16515890193334723" @("eyes open", "eyes closed") $s3
$s5 = New-TaskOrderSheet "GNG_TO-16515890193647213" @("go_stims-16515890193334723.csv", "GNG_stims-16515890193490965.csv", "go_stims-16515890193490965.csv", "GNG_stims-16515890193647213.csv") $s4

# Remove the original sheets now that their replacements are in place.
for ($i = 1; $i -le $origCount; $i++) {
    $wb.Worksheets.Item(1).Delete() | Out-Null
}

# Re-fetch by position: handles captured before the deletions above can go
# stale once the sheet collection is mutated, so look the first tab up fresh.
$wb.Worksheets.Item(1).Activate()
